$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '49.646.04'
$ws.Range('E2').Value = '  -1.02%  '
$ws.Range('D3').Value = '2.640.51'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '113.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.92%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '324.62'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.545'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.13%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.83'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0814'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E13').Value = '  +1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.33'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.29%  '
$ws.Range('D15').Value = '3.058.29'
$ws.Range('E15').Value = '  -0.11%  '
$ws.Range('D16').Value = '2.639.68'
$ws.Range('E16').Value = '  +0.14%  '
$ws.Range('E17').Value = '  -2.43%  '
$ws.Range('D18').Value = '49.597.05'
$ws.Range('E18').Value = '  -0.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.99'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.92%  '
$ws.Range('E21').Value = '  -1.80%  '
$ws.Range('D22').Value = '0.0₃0949'
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '270.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.99'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.64%  '
$ws.Range('E25').Value = '  -2.71%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.35'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.24%  '
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.34'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.94%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.23'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.27'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.75%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.139'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.59'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.63%  '
$ws.Range('E33').Value = '  +0.75%  '
$ws.Range('E34').Value = '  +2.29%  '
$ws.Range('E35').Value = '  -0.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.07'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.96'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.11%  '
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.13'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '127.04'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.74%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.111'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0327'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +3.92%  '
$ws.Range('E44').Value = '  -3.45%  '
$ws.Range('D45').Value = '2.064.71'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.24'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.29%  '
$ws.Range('B47').Value = 'Stacks'
$ws.Range('C47').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.13'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.43%  '
$ws.Range('E48').Value = '  -6.94%  '
$ws.Range('E49').Value = '  -1.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '59.32'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.27%  '
$ws.Range('E51').Value = '  -3.12%  '
